$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the number-format styling from the row above (row 5) so the new
# date/boolean cells reuse the existing style (s="1") instead of creating
# a brand-new numFmt entry in styles.xml.
$ws.Range("A5").Copy()
$ws.Range("A6").PasteSpecial(-4122)

$ws.Range("G5").Copy()
$ws.Range("G6").PasteSpecial(-4122)

# Now populate the new row's values.
$ws.Range("A6").Value = 42647.680659722224
$ws.Range("B6").Value = $true
$ws.Range("C6").Value = 9941.89
$ws.Range("D6").Value = 9766.58
$ws.Range("E6").Value = 18.12
$ws.Range("F6").Value = 18.77
$ws.Range("G6").Value = $false
$ws.Range("H6").Value = 3.59
$ws.Range("I6").Value = $false

$excel.CutCopyMode = $false
